$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final")

# Append the four new teams/rows to the bottom of the "Final" standings table.
$ws.Range("A10").Value = "I Love Bong Pitts"
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = 8
$ws.Range("D10").Value = 1776

$ws.Range("A11").Value = "London Calling"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 1678

$ws.Range("A12").Value = "Dak White"
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 1620

$ws.Range("A13").Value = "Finnegan's Fantastic Team"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 12
$ws.Range("D13").Value = 1334

# Column A now needs to fit the longest team name ("Finnegan's Fantastic Team"),
# same as the already-present column A on the "Records" sheet.
$ws.Columns.Item(1).ColumnWidth = 21.6

# Move the active selection to just past the new data, matching the saved view.
$ws.Range("D14").Select() | Out-Null
